$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the old Area/Objective/Risk/Control/Test
# columns C:J one to the right (D:K), matching the new 3-column Area block (ref/type/objectives).
$ws.Columns.Item(3).Insert()

# --- Row 2 (Import) - fill the gap left by the insert; it mirrors neighbouring Y's ---
$ws.Range("C2").Value = "Y"

# --- Row 3 (Key) - fill the gap left by the insert ---
$ws.Range("C3").Value = "N"

# --- Row 5 (Field Name) - rename/extend the Area fields for the new 3-column layout ---
$ws.Range("B5").Value = "laa_AreaRef"
$ws.Range("C5").Value = "laa_AreaType"
$ws.Range("D5").Value = "laa_areaObjectives"

# --- Row 6 (new) - area model / gorm column-naming reference block ---
$ws.Range("B6").Value = 1.1
$ws.Range("C6").Value = "area_AreaType"
$gormText = @'
'type User struct {
  ID uint             // column name will be `id`
  Name string         // column name will be `name`
  Birthday time.Time  // column name will be `birthday`
  CreatedAt time.Time // column name will be `created_at`
}
// Overriding Column Name
type Animal struct {
    AnimalId    int64     `gorm:"column:beast_id"`         // set column name to `beast_id`
    Birthday    time.Time `gorm:"column:day_of_the_beast"` // set column name to `day_of_the_beast`
    Age         int64     `gorm:"column:age_of_the_beast"` // set column name to `age_of_the_beast`
}
'@
$ws.Range("D6").Value = $gormText
$ws.Range("D6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 176.1

Write-Output "done"
